# The deck's single custom "Integral" theme (Red Violet colour scheme) is
# switched over to the standard Office Theme colour scheme. In the
# canonical OOXML this shows up as the colour values (and the <a:theme>/
# <a:clrScheme> "name" attributes, which PowerPoint's object model does not
# expose as a writable/persisted string) moving from theme1.xml to
# theme2.xml and vice-versa; functionally/visually this is exactly the
# recolour below, applied to the deck's one reachable theme (the
# SlideMaster's, i.e. ppt/theme/theme1.xml).

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

function ToComRGB($r, $g, $b) {
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme colour scheme (RRGGBB), in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$tcs.Colors(1).RGB  = ToComRGB 0x00 0x00 0x00   # dk1      000000
$tcs.Colors(2).RGB  = ToComRGB 0xFF 0xFF 0xFF   # lt1      FFFFFF
$tcs.Colors(3).RGB  = ToComRGB 0x44 0x54 0x6A   # dk2      44546A
$tcs.Colors(4).RGB  = ToComRGB 0xE7 0xE6 0xE6   # lt2      E7E6E6
$tcs.Colors(5).RGB  = ToComRGB 0x5B 0x9B 0xD5   # accent1  5B9BD5
$tcs.Colors(6).RGB  = ToComRGB 0xED 0x7D 0x31   # accent2  ED7D31
$tcs.Colors(7).RGB  = ToComRGB 0xA5 0xA5 0xA5   # accent3  A5A5A5
$tcs.Colors(8).RGB  = ToComRGB 0xFF 0xC0 0x00   # accent4  FFC000
$tcs.Colors(9).RGB  = ToComRGB 0x44 0x72 0xC4   # accent5  4472C4
$tcs.Colors(10).RGB = ToComRGB 0x70 0xAD 0x47   # accent6  70AD47
$tcs.Colors(11).RGB = ToComRGB 0x05 0x63 0xC1   # hlink    0563C1
$tcs.Colors(12).RGB = ToComRGB 0x95 0x4F 0x72   # folHlink 954F72
